# BIS-1002: Fixed XLS export tests
# Adds a new "Internal Assignment" column (O) to the property-type table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new column, styled like the other header cells (N4)
# but bold + font size 12.
$ws.Range("O4").Value = "Internal Assignment"
$ws.Range("O4").Font.Bold = $true
$ws.Range("O4").Font.Size = 12
$ws.Range("O4").Font.Name = "Calibri"
$ws.Range("O4").Font.Color = 0

# Data cells for the property rows.
$ws.Range("O5").Value = "TRUE"
$ws.Range("O6").Value = "FALSE"
$ws.Range("O7").Value = "FALSE"

$ws.Range("O5:O7").Font.Name = "Calibri"
$ws.Range("O5:O7").Font.Size = 11
$ws.Range("O5:O7").Font.Bold = $false

$ws.Range("O6").Select()
